# UsersAPI.xlsx edit:
#  - reorder worksheet tabs so "PostUsers" comes before "PostUsersExistOneField"
#  - make "PostUsers" the active tab
#  - update sample data (first names in PostUsers, and first name / contact /
#    email / patient-id in PostUsersExistOneField)
#  - refresh the UI selections / scroll position left behind by Excel on the
#    four sheets

$wb = $excel.ActiveWorkbook

# --- 1. Move "PostUsers" in front of "PostUsersExistOneField" -------------
$postUsers = $wb.Worksheets.Item("PostUsers")
$postUsers.Move($wb.Worksheets.Item(1))

# --- 2. Update PostUsersExistOneField sample data --------------------------
$existOneField = $wb.Worksheets.Item("PostUsersExistOneField")
$existOneField.Range("A2").Value = "Sneha"
$existOneField.Range("L2").Value = "PT2630"
$existOneField.Range("D3").Value = 1234067089
$existOneField.Range("L3").Value = "PT5248"
$existOneField.Range("E4").Value = "abcs@xyz.com"
$existOneField.Range("L4").Value = "PT6346"

# --- 3. Update PostUsers sample data ---------------------------------------
$postUsers = $wb.Worksheets.Item("PostUsers")
$postUsers.Range("A2").Value = "Kirti"
$postUsers.Range("A3").Value = "Shane"
$postUsers.Range("A4").Value = "Shaun"
$postUsers.Range("A5").Value = "Duncan"
$postUsers.Range("A6").Value = "Priya"

# --- 4. Fix up the leftover UI state (selection / scroll / active tab) -----
# Touch the non-active sheets first, then finish on "PostUsers" so it ends up
# as the active/selected tab.
$existOneField.Range("D24").Select()

$userType = $wb.Worksheets.Item("PostUsersUserType")
$userType.Range("H3").Select()

$missingFields = $wb.Worksheets.Item("PostUsersMissingFields")
$missingFields.Range("A2").Select()

$postUsers.Range("L1").Select()
$postUsers.Activate()
